# "Generate Report for Archive"
#
# The localization-status report is regenerated: every cell that showed the
# "Ready for handoff" status now reads "In Translation" (Overview!E2:F2,
# zh-cn!C2, de-de!C2). Excel auto-fits the status columns to the new
# (shorter) text, which narrows them from ~17.22 chars to ~13.41 chars.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# Target column width (in Excel's stored "character" width units) that a
# real AutoFit settles on for the new, shorter status text.
$targetStoredWidth = 13.4101845877511
# This host's Range.ColumnWidth setter stores (input + 5/6); back it out so
# the saved width lands as close as possible to the real AutoFit result.
$newColumnWidth = $targetStoredWidth - (5 / 6)

# --- Overview sheet: status lives in columns E (zh-cn) and F (de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth

# --- Per-locale sheets: status lives in column C ---
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("C2").Value = $newStatus
    $ws.Columns.Item(3).ColumnWidth = $newColumnWidth
}
